$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(1).ColumnWidth = 28.877604166666668
$ws.Columns.Item(2).ColumnWidth = 34.022135416666664
$ws.Columns.Item(3).ColumnWidth = 33.307291666666664
$ws.Columns.Item(4).ColumnWidth = 32.307291666666664
$ws.Columns.Item(5).ColumnWidth = 13.736979166666666
Write-Host "done"
